$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 23420
$ws.Range("I57").Value = 30000
$ws.Range("J57").Value = 21775
$ws.Range("K57").Value = 90000
$ws.Range("L57").Value = 65325
$ws.Range("M57").Value = -89501
$ws.Range("N57").Value = -66323
$ws.Range("H75").Value = 24552.5
$ws.Range("J75").Value = 29063
$ws.Range("L75").Value = 29063
$ws.Range("N75").Value = -30935
$ws.Range("H78").Value = 24552.5
$ws.Range("J78").Value = 29063
$ws.Range("L78").Value = 87189
$ws.Range("N78").Value = -96549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2405.2654
$ws.Range("I74").Value = 2689.8572
$ws.Range("J74").Value = 1693.7858
$ws.Range("K74").Value = 2689.8572
$ws.Range("L74").Value = 1693.7858
$ws.Range("M74").Value = -1815.8572
$ws.Range("N74").Value = -3441.7858
$ws.Range("H77").Value = 2405.2654
$ws.Range("I77").Value = 2689.8572
$ws.Range("J77").Value = 1693.7858
$ws.Range("K77").Value = 13449.286
$ws.Range("L77").Value = 8468.929
$ws.Range("M77").Value = -9081.286
$ws.Range("N77").Value = -17204.929
$ws.Range("H132").Value = 1527.2712
$ws.Range("I132").Value = 975.6981
$ws.Range("J132").Value = 6399.5
$ws.Range("K132").Value = 2927.0943
$ws.Range("L132").Value = 19198.5
$ws.Range("M132").Value = -397.0942999999997
$ws.Range("N132").Value = -24258.5
$ws.Range("H139").Value = 44200
$ws.Range("J139").Value = 44200
$ws.Range("L139").Value = 44200
$ws.Range("N139").Value = -54480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 41380
$ws.Range("J138").Value = 41380
$ws.Range("L138").Value = 41380
$ws.Range("N138").Value = -51660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11908002
$ws.Range("I31").Value = 1665.6818
$ws.Range("J31").Value = 25004972
$ws.Range("K31").Value = 1665.6818
$ws.Range("L31").Value = 25004972
$ws.Range("M31").Value = -1370.6818
$ws.Range("N31").Value = -25005562
$ws.Range("H34").Value = 11908002
$ws.Range("I34").Value = 1665.6818
$ws.Range("J34").Value = 25004972
$ws.Range("K34").Value = 1665.6818
$ws.Range("L34").Value = 25004972
$ws.Range("M34").Value = -1463.6818
$ws.Range("N34").Value = -25005376
$ws.Range("H58").Value = 1450.75
$ws.Range("I58").Value = 1396.8353
$ws.Range("J58").Value = 1756.2667
$ws.Range("K58").Value = 1396.8353
$ws.Range("L58").Value = 1756.2667
$ws.Range("M58").Value = -1193.8353
$ws.Range("N58").Value = -2162.2667
$ws.Range("H134").Value = 1361.5538
$ws.Range("I134").Value = 795.8333
$ws.Range("J134").Value = 2958.8823
$ws.Range("K134").Value = 2387.4999
$ws.Range("L134").Value = 8876.6469
$ws.Range("M134").Value = 147.5001000000002
$ws.Range("N134").Value = -13946.6469
$ws.Range("H136").Value = 1450.75
$ws.Range("I136").Value = 1396.8353
$ws.Range("J136").Value = 1756.2667
$ws.Range("K136").Value = 4190.5059
$ws.Range("L136").Value = 5268.800099999999
$ws.Range("M136").Value = -1640.5059
$ws.Range("N136").Value = -10368.8001
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H140").Value = 140130
$ws.Range("J140").Value = 140130
$ws.Range("L140").Value = 140130
$ws.Range("N140").Value = -150490
$ws.Range("H141").Value = 35000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 35000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = ""
$ws.Range("M141").Value = 35000
$ws.Range("N141").Value = -45360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 3310.9412
$ws.Range("I21").Value = 629.3333
$ws.Range("J21").Value = 3885.5715
$ws.Range("K21").Value = 1887.9999
$ws.Range("L21").Value = 11656.7145
$ws.Range("M21").Value = -1714.9999
$ws.Range("N21").Value = -12002.7145
$ws.Range("H26").Value = 7537.35
$ws.Range("I26").Value = 12154.8
$ws.Range("J26").Value = 2919.9
$ws.Range("K26").Value = 36464.39999999999
$ws.Range("L26").Value = 8759.700000000001
$ws.Range("M26").Value = -36176.39999999999
$ws.Range("N26").Value = -9335.700000000001
$ws.Range("H37").Value = 333466660
$ws.Range("J37").Value = 333466660
$ws.Range("L37").Value = 1000399980
$ws.Range("N37").Value = -1000400204
$ws.Range("H113").Value = 616.807
$ws.Range("I113").Value = 520.43475
$ws.Range("J113").Value = 1019.8182
$ws.Range("K113").Value = 1561.30425
$ws.Range("L113").Value = 3059.4546
$ws.Range("M113").Value = 608.6957499999999
$ws.Range("N113").Value = -7399.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 37508.5
$ws.Range("I25").Value = 6008
$ws.Range("J25").Value = 69009
$ws.Range("K25").Value = 6008
$ws.Range("L25").Value = 69009
$ws.Range("M25").Value = -5479
$ws.Range("N25").Value = -70067
$ws.Range("H35").Value = 34994
$ws.Range("J35").Value = 34994
$ws.Range("L35").Value = 34994
$ws.Range("N35").Value = -35590
$ws.Range("H39").Value = 17862.25
$ws.Range("J39").Value = 17862.25
$ws.Range("L39").Value = 17862.25
$ws.Range("N39").Value = -18926.25
$ws.Range("H41").Value = 11271.667
$ws.Range("J41").Value = 21560.5
$ws.Range("L41").Value = 21560.5
$ws.Range("N41").Value = -22270.5
$ws.Range("H96").Value = 26660.334
$ws.Range("J96").Value = 26660.334
$ws.Range("L96").Value = 26660.334
$ws.Range("N96").Value = -32152.334
$ws.Range("H100").Value = 39380
$ws.Range("J100").Value = 39380
$ws.Range("L100").Value = 39380
$ws.Range("N100").Value = -41544
$ws.Range("H102").Value = 1599.7
$ws.Range("I102").Value = 1108.8462
$ws.Range("J102").Value = 2511.2856
$ws.Range("K102").Value = 1108.8462
$ws.Range("L102").Value = 2511.2856
$ws.Range("M102").Value = 513.1538
$ws.Range("N102").Value = -5755.2856
$ws.Range("H140").Value = 38671.613
$ws.Range("J140").Value = 38671.613
$ws.Range("L140").Value = 38671.613
$ws.Range("N140").Value = -49031.613

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 15004
$ws.Range("J13").Value = 15004
$ws.Range("L13").Value = 15004
$ws.Range("N13").Value = -15284
$ws.Range("H25").Value = 10897.5
$ws.Range("I25").Value = 6797
$ws.Range("J25").Value = 14998
$ws.Range("K25").Value = 6797
$ws.Range("L25").Value = 14998
$ws.Range("M25").Value = -6567
$ws.Range("N25").Value = -15458
$ws.Range("H42").Value = 34999
$ws.Range("J42").Value = 34999
$ws.Range("L42").Value = 34999
$ws.Range("N42").Value = -36125
$ws.Range("H49").Value = 34999
$ws.Range("J49").Value = 34999
$ws.Range("L49").Value = 34999
$ws.Range("N49").Value = -35293
$ws.Range("H74").Value = 33198.375
$ws.Range("I74").Value = 7098.5
$ws.Range("J74").Value = 41898.332
$ws.Range("K74").Value = 7098.5
$ws.Range("L74").Value = 41898.332
$ws.Range("M74").Value = -6100.5
$ws.Range("N74").Value = -43894.332
$ws.Range("H76").Value = 29799.334
$ws.Range("J76").Value = 29799.334
$ws.Range("L76").Value = 29799.334
$ws.Range("N76").Value = -30475.334
$ws.Range("H77").Value = 33198.375
$ws.Range("I77").Value = 7098.5
$ws.Range("J77").Value = 41898.332
$ws.Range("K77").Value = 21295.5
$ws.Range("L77").Value = 125694.996
$ws.Range("M77").Value = -16303.5
$ws.Range("N77").Value = -135678.996
$ws.Range("H79").Value = 29799.334
$ws.Range("J79").Value = 29799.334
$ws.Range("L79").Value = 29799.334
$ws.Range("N79").Value = -32139.334
$ws.Range("H139").Value = 40910
$ws.Range("J139").Value = 40910
$ws.Range("L139").Value = 40910
$ws.Range("N139").Value = -51190
$ws.Range("H140").Value = 70295.22
$ws.Range("J140").Value = 70295.22
$ws.Range("L140").Value = 70295.22
$ws.Range("N140").Value = -80655.22
$ws.Range("H141").Value = 41723
$ws.Range("J141").Value = 41723
$ws.Range("L141").Value = 41723
$ws.Range("N141").Value = -52083

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 49352.668
$ws.Range("I37").Value = 8000
$ws.Range("J37").Value = 70029
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 70029
$ws.Range("M37").Value = -7797
$ws.Range("N37").Value = -70435
$ws.Range("H136").Value = 2872.85
$ws.Range("I136").Value = 1039.1904
$ws.Range("J136").Value = 4899.5264
$ws.Range("K136").Value = 3117.5712
$ws.Range("L136").Value = 14698.5792
$ws.Range("M136").Value = -567.5711999999999
$ws.Range("N136").Value = -19798.5792
$ws.Range("H138").Value = 44568
$ws.Range("J138").Value = 44568
$ws.Range("L138").Value = 44568
$ws.Range("N138").Value = -54848
$ws.Range("H139").Value = 39099.78
$ws.Range("I139").Value = 35325
$ws.Range("J139").Value = 39459.285
$ws.Range("K139").Value = 35325
$ws.Range("L139").Value = 39459.285
$ws.Range("M139").Value = -30185
$ws.Range("N139").Value = -49739.285
$ws.Range("H140").Value = 40262.332
$ws.Range("J140").Value = 40262.332
$ws.Range("L140").Value = 40262.332
$ws.Range("N140").Value = -50622.332
$ws.Range("H141").Value = 42433.2
$ws.Range("J141").Value = 42433.2
$ws.Range("L141").Value = 42433.2
$ws.Range("N141").Value = -52793.2
